# This workbook logs daily "Palta" (avocado) price records. The edit
# inserts one new daily record before the existing row 315, which pushes
# row 315 and everything below it down by one row (the former last row,
# 421, becomes the new row 422). This matches the canonical diff, where
# the dimension grows from A1:T421 to A1:T422 and every row from 315 on
# is shifted down by one, with a brand-new record appearing at (new) row
# 315.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 315; rows 315..421 shift to 316..422,
# and formatting (e.g. the date style on column D) is inherited from the
# row above, matching the original file's styling.
$ws.Rows("315:315").Insert()

# Columns A-L of the new record are identical to the record that used to
# sit at row 315 (now shifted to row 316): same market, region, codreg,
# type, product, category, variety ("Hass") and quality ("Primera").
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(315, $c).Value = $ws.Cells.Item(316, $c).Value2
}

# Columns that differ for the new record: date, volume, min/max/avg
# price, commercialization unit, origin, $/kg and kg per unit.
$ws.Cells.Item(315, 4).Value  = 44588
$ws.Cells.Item(315, 13).Value = 200
$ws.Cells.Item(315, 14).Value = 2600
$ws.Cells.Item(315, 15).Value = 2800
$ws.Cells.Item(315, 16).Value = 2700
$ws.Cells.Item(315, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(315, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(315, 19).Value = 2700
$ws.Cells.Item(315, 20).Value = 1
